$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -6.059299999999999
$ws.Range("D21").Value = -7.639000000000001
$ws.Range("D23").Value = -6.903499999999992
$ws.Range("D25").Value = -8.342399999999998
